$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "70.886.24"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.32%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.803.59"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.66%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "705.13"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.16%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "170.04"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -2.37%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.803.52"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -1.60%  "

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.02%  "

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.78%  "

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -1.39%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.63"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +5.77%  "

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.48%  "

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -3.44%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.78"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -1.94%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.446.57"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -1.56%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.772.90"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -2.45%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "70.869.09"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.39%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "17.36"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -1.62%  "

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.34%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.09"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -1.81%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "499.18"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.29%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.64"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -2.84%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.722"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.07%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "83.97"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -1.14%  "

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -5.25%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.954.04"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -1.32%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.01"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -1.88%  "

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -4.10%  "

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.02%  "

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -6.49%  "

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -4.54%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.25"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -1.69%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.30"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -3.62%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "28.96"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -2.47%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.173"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -5.45%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.00"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.01%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.769.47"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -1.28%  "

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -2.25%  "

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -3.63%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.34"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -2.11%  "

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -1.36%  "

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -3.52%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.23"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -5.86%  "

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.16%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "166.76"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +1.88%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.000313"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.56%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "49.03"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.86%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "418.39"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.32%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.57"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.80%  "

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -3.06%  "
